$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17, shifting the existing data (old rows 17-128)
# down to rows 18-129. This matches Excel's default "insert row" behaviour,
# which also shifts cell formatting (e.g. the date style on column D) down
# with the rest of the row.
$ws.Rows(17).Insert(-4121)

# Populate the freshly inserted row 17 with the new weekly data point.
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44635
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112044
$ws.Range("G17").Value = "Perejil"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 2400
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 2750
$ws.Range("N17").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O17").Value = "Provincia del Elqu" + [char]237
$ws.Range("P17").Value = 1833
$ws.Range("Q17").Value = 1.5
$ws.Range("R17").Value = "Hortaliza"
